# Add a "科室" (Department) column to the RegionMeetings export template.
# The new column is inserted before the existing "${record.province}"
# column (old column E), shifting every following column one place to
# the right - matching the upstream commit "Add missing department to
# export data".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E; Excel shifts D:N -> E:O automatically,
# carrying along values, styles and the shared-string table.
$ws.Columns("E:E").Insert()

# Populate the new header (row 1) and template placeholder (row 2) cells.
$ws.Range("E1").Value = "科室"
$ws.Range("E2").Value = '${record.department}'

# The jx: template markers recorded in the header comments hard-code the
# last cell of the generated area/each block as text, so they need to be
# updated by hand to reflect the newly widened N2 -> O2 range.
[void]$ws.Range("A1").Comment.Text("Author:`njx:area(lastCell=`"O2`")")
[void]$ws.Range("A2").Comment.Text("Author:`njx:each(items=`"data`" var=`"record`" lastCell=`"O2`")")

# Refresh the view: select E10 (clears the stale topLeftCell/selection
# left over from before the column insert).
[void]$ws.Range("E10").Select()
